$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Insert 5 new rows above the old row 9 ("BASIC ROM (reserved)"), pushing the
# old rows 9/10 down to 14/15. Excel carries the surviving rows' data,
# formulas and formatting down automatically.
# ---------------------------------------------------------------------------
$ws.Rows("9:13").Insert()

# ---------------------------------------------------------------------------
# New row 9: ASM (entry)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value2 = 2049
$ws.Range("C9").Formula = "=_xlfn.CONCAT(""$"", DEC2HEX(B9,4))"
$ws.Range("D9").Value2 = 258
$ws.Range("E9").Formula = "=B9+D9-1"
$ws.Range("F9").Formula = "=_xlfn.CONCAT(""$"",DEC2HEX(E9,4))"
$ws.Range("G9").Formula = "=IF(E9>=B10,""yes"",""no"")"
$ws.Range("H9").Formula = "=IF(B9<=16383,0,IF(B9<=32767,1,IF(B9<=49151,2,3)))"
$ws.Range("I9").Formula = "=IF((B9+D9-1) < ((H9+1)*16384),""yes"", ""no"")"

# ---------------------------------------------------------------------------
# New row 10: sprite data
# ---------------------------------------------------------------------------
$ws.Range("B10").Value2 = 11904
$ws.Range("C10").Formula = "=_xlfn.CONCAT(""$"", DEC2HEX(B10,4))"
$ws.Range("D10").Formula = "=6*64"
$ws.Range("E10").Formula = "=B10+D10-1"
$ws.Range("F10").Formula = "=_xlfn.CONCAT(""$"",DEC2HEX(E10,4))"
$ws.Range("G10").Formula = "=IF(E10>=B11,""yes"",""no"")"
$ws.Range("H10").Formula = "=IF(B10<=16383,0,IF(B10<=32767,1,IF(B10<=49151,2,3)))"
$ws.Range("I10").Formula = "=IF((B10+D10-1) < ((H10+1)*16384),""yes"", ""no"")"

# ---------------------------------------------------------------------------
# New row 11: CHARDATA
# ---------------------------------------------------------------------------
$ws.Range("B11").Value2 = 12288
$ws.Range("C11").Formula = "=_xlfn.CONCAT(""$"", DEC2HEX(B11,4))"
$ws.Range("D11").Value2 = 128
$ws.Range("E11").Formula = "=B11+D11-1"
$ws.Range("F11").Formula = "=_xlfn.CONCAT(""$"",DEC2HEX(E11,4))"
$ws.Range("G11").Formula = "=IF(E11>=B12,""yes"",""no"")"
$ws.Range("H11").Formula = "=IF(B11<=16383,0,IF(B11<=32767,1,IF(B11<=49151,2,3)))"
$ws.Range("I11").Formula = "=IF((B11+D11-1) < ((H11+1)*16384),""yes"", ""no"")"

# ---------------------------------------------------------------------------
# New row 12: _color_data
# ---------------------------------------------------------------------------
$ws.Range("B12").Value2 = 38912
$ws.Range("C12").Formula = "=_xlfn.CONCAT(""$"", DEC2HEX(B12,4))"
$ws.Range("D12").Value2 = 1024
$ws.Range("E12").Formula = "=B12+D12-1"
$ws.Range("F12").Formula = "=_xlfn.CONCAT(""$"",DEC2HEX(E12,4))"
$ws.Range("G12").Formula = "=IF(E12>=B13,""yes"",""no"")"
$ws.Range("H12").Formula = "=IF(B12<=16383,0,IF(B12<=32767,1,IF(B12<=49151,2,3)))"
$ws.Range("I12").Formula = "=IF((B12+D12-1) < ((H12+1)*16384),""yes"", ""no"")"

# ---------------------------------------------------------------------------
# New row 13: _screen_data
# ---------------------------------------------------------------------------
$ws.Range("B13").Value2 = 39936
$ws.Range("C13").Formula = "=_xlfn.CONCAT(""$"", DEC2HEX(B13,4))"
$ws.Range("D13").Value2 = 1024
$ws.Range("E13").Formula = "=B13+D13-1"
$ws.Range("F13").Formula = "=_xlfn.CONCAT(""$"",DEC2HEX(E13,4))"
$ws.Range("G13").Formula = "=IF(E13>=B14,""yes"",""no"")"
$ws.Range("H13").Formula = "=IF(B13<=16383,0,IF(B13<=32767,1,IF(B13<=49151,2,3)))"
$ws.Range("I13").Formula = "=IF((B13+D13-1) < ((H13+1)*16384),""yes"", ""no"")"

# ---------------------------------------------------------------------------
# Row 15 (old row 10, "I/O registers (reserved)") previously had a blank G
# cell; the fill handle now reaches one row further than the data, so G15
# gets a real formula referencing the (blank) B16.
# ---------------------------------------------------------------------------
$ws.Range("G15").Formula = "=IF(E15>=B16,""yes"",""no"")"

# ---------------------------------------------------------------------------
# Description labels (shared strings) - set in this specific order so the
# new shared-string pool matches the source ordering (CHARDATA, _screen_data,
# _color_data, sprite data, ASM (entry)).
# ---------------------------------------------------------------------------
$ws.Range("A11").Value2 = "CHARDATA"
$ws.Range("A13").Value2 = "_screen_data"
$ws.Range("A12").Value2 = "_color_data"
$ws.Range("A10").Value2 = "sprite data"
$ws.Range("A9").Value2 = "ASM (entry)"

# ---------------------------------------------------------------------------
# Formatting: rows 9-15 columns G/H/I did not carry a highlight fill in the
# source; make sure that's the case here too (matches the "no fill" cellXfs
# added to styles.xml for these columns).
# ---------------------------------------------------------------------------
$ws.Range("G9:I15").Interior.ColorIndex = -4142

# Row 9 description/value cells use the bold header font (like row 1); keep
# the numeric columns using the regular "highlighted" look used elsewhere.
$ws.Range("A9:B9").Font.Bold = $true
$ws.Range("E9").Font.Bold = $true

# ---------------------------------------------------------------------------
# View / selection bookkeeping
# ---------------------------------------------------------------------------
$ws.Range("A11").Select()
